$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concentrations")

# Insert a new column D ("Distribution") into the Concentrations sheet, in front
# of the existing Mean / Standard deviation / Units columns.
$ws.Columns.Item(4).Insert()

# Header for the new column.
$ws.Cells.Item(1, 4).Value = "Distribution"

# Fill the new column with the distribution type used for every row (all
# initial concentrations in this fixture are normally distributed).
$lastRow = 7
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "normal"
}

# Rename the sheet to clarify the semantics of the values it holds. Excel
# worksheet names are limited to 31 characters, so use the longest valid
# prefix of the intended name.
$ws.Name = "Distributions of initial concen"

# Make this sheet the active / selected tab, with the last-used cell
# selection set to M14 (matching the workbook's recorded view state).
$ws.Activate()
$ws.Range("M14").Select()
